# nay - Movie API, Screen Fix
#
# 1) API一覧 (sheet2): insert a new row before row 18 for the new
#    "(管理用)同タイトルムービーチェック" movie-title-check API, copying the
#    formatting of the row above (row 17) into the newly inserted row, then
#    column widths are widened for columns A and E.
# 2) 機能画面一覧 (sheet1): active selection changes, and the active sheet
#    moves from sheet1 to sheet2 (API一覧).

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# --- API一覧: insert the new "movie title check" API row above row 18 ---
$ws2.Range("A18:I18").Insert(-4121)

# Copy the formatting from the (unshifted) row above the insertion point
# (row 17) down onto the freshly inserted row 18 so the new row matches the
# surrounding table styling.
for ($col = 1; $col -le 9; $col++) {
    $ws2.Cells.Item(17, $col).Copy()
    $ws2.Cells.Item(18, $col).PasteSpecial(-4122)
}

$ws2.Cells.Item(18, 2).Value = "/admin/movie/title/{title}"
$ws2.Cells.Item(18, 1).Value = "（管理用）同タイトルムービーチェック"
$ws2.Cells.Item(18, 3).Value = "GET"
$ws2.Cells.Item(18, 5).Value = "ムービー作成（管理用）"

# Widen column A and column E on the API一覧 sheet.
$ws2.Columns.Item(1).ColumnWidth = 40.21875
$ws2.Columns.Item(5).ColumnWidth = 44.5546875

# --- 機能画面一覧: update the remembered selection on that sheet ---
$ws1.Activate()
$ws1.Range("A15:J15").Select()

# --- API一覧 becomes the active tab, with its own remembered selection ---
$ws2.Activate()
$ws2.Range("E21").Select()
